$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "305.12"
Set-TextValue $ws.Range("E2") "6.70%"
Set-TextValue $ws.Range("D3") "31.94"
Set-TextValue $ws.Range("E3") "8.80%"
Set-TextValue $ws.Range("D4") "5.279"
Set-TextValue $ws.Range("E4") "3.85%"
Set-TextValue $ws.Range("D5") "0.07511"
Set-TextValue $ws.Range("E5") "11.77%"
Set-TextValue $ws.Range("D6") "7.844"
Set-TextValue $ws.Range("D7") "3.762"
Set-TextValue $ws.Range("E7") "9.49%"
Set-TextValue $ws.Range("D8") "1.471"
Set-TextValue $ws.Range("E8") "6.54%"
Set-TextValue $ws.Range("D9") "0.9160"
Set-TextValue $ws.Range("E9") "2.38%"
Set-TextValue $ws.Range("D10") "0.01753"
Set-TextValue $ws.Range("E10") "2,614.46%"
Set-TextValue $ws.Range("D11") "0.1693"
Set-TextValue $ws.Range("E11") "7.40%"
Set-TextValue $ws.Range("D12") "0.07828"
Set-TextValue $ws.Range("E12") "14.52%"
Set-TextValue $ws.Range("D13") "0.08045"
Set-TextValue $ws.Range("E13") "5.47%"
Set-TextValue $ws.Range("D14") "0.03015"
Set-TextValue $ws.Range("E14") "3.13%"
Set-TextValue $ws.Range("D15") "0.09892"
Set-TextValue $ws.Range("E15") "9.97%"
Set-TextValue $ws.Range("D16") "0.001493"
Set-TextValue $ws.Range("E16") "-5.61%"
Set-TextValue $ws.Range("D17") "0.04556"
Set-TextValue $ws.Range("E17") "1.89%"
Set-TextValue $ws.Range("D18") "0.006295"
Set-TextValue $ws.Range("E18") "1.32%"
Set-TextValue $ws.Range("D19") "3.483"
Set-TextValue $ws.Range("E19") "0.96%"
Set-TextValue $ws.Range("D20") "2.231"
Set-TextValue $ws.Range("E20") "0.01%"
Set-TextValue $ws.Range("D21") "0.3301"
Set-TextValue $ws.Range("E21") "3.02%"
Set-TextValue $ws.Range("D22") "0.1343"
Set-TextValue $ws.Range("E22") "1.82%"
Set-TextValue $ws.Range("D23") "4.476"
Set-TextValue $ws.Range("E23") "8.96%"
Set-TextValue $ws.Range("D24") "0.1637"
Set-TextValue $ws.Range("E24") "3.61%"
Set-TextValue $ws.Range("D25") "0.001216"
Set-TextValue $ws.Range("E25") "1.07%"
Set-TextValue $ws.Range("D26") "0.004432"
Set-TextValue $ws.Range("E26") "1.36%"
Set-TextValue $ws.Range("D27") "0.0001398"
Set-TextValue $ws.Range("E27") "19.62%"
Set-TextValue $ws.Range("D28") "0.0001739"
Set-TextValue $ws.Range("E28") "7.58%"
Set-TextValue $ws.Range("D40") "0.04526"
Set-TextValue $ws.Range("E40") "6.64%"
Set-TextValue $ws.Range("D41") "0.007144"
Set-TextValue $ws.Range("E41") "5.19%"
Set-TextValue $ws.Range("D42") "0.1345"
Set-TextValue $ws.Range("E42") "8.69%"
Set-TextValue $ws.Range("D43") "0.002247"
Set-TextValue $ws.Range("E43") "2.71%"
Set-TextValue $ws.Range("D44") "0.01399"
Set-TextValue $ws.Range("E44") "22.32%"
Set-TextValue $ws.Range("E45") "8.29%"
Set-TextValue $ws.Range("D46") "1.870"
Set-TextValue $ws.Range("E46") "-2.85%"
Set-TextValue $ws.Range("D47") "0.01299"
Set-TextValue $ws.Range("E47") "-13.47%"
